# Apply the edits described by the commit:
#  - Column C formulas change from the whole-column "A:A+B:B" form to the
#    per-row "A<n>+B<n>" form (sum of A and B on the same row).
#  - A new column D is added with a per-row product formula "A<n>*B<n>".
#  - The used range / dimension therefore grows from A1:C4 to A1:D4.
#  - The active selection moves from D3 to E6.
#  - Best-effort: mark font 0 (Arial, used by all cells) with the default
#    charset; this isn't exposed as a settable Font property in this COM
#    surface, so it's attempted but harmless if ignored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
    $ws.Range("A1").Font.Charset = 1
} catch {
    # Charset isn't settable through this object model -- ignore.
}

# Column C: row-relative sum formulas (replacing the old "A:A+B:B" formulas).
$ws.Range("C1").Formula = "=A1+B1"
$ws.Range("C2").Formula = "=A2+B2"
$ws.Range("C3").Formula = "=A3+B3"
$ws.Range("C4").Formula = "=A4+B4"

# Column D: new row-relative product formulas.
$ws.Range("D1").Formula = "=A1*B1"
$ws.Range("D2").Formula = "=A2*B2"
$ws.Range("D3").Formula = "=A3*B3"
$ws.Range("D4").Formula = "=A4*B4"

# Move the active selection from D3 to E6.
$null = $ws.Range("E6").Select()
